$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 8 (Serbian Super League, Italian Serie B)
# so existing rows 8-12 shift down to rows 10-14.
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

# Row 2: odds updates
$ws.Cells.Item(2, 6).Value = 2.82
$ws.Cells.Item(2, 7).Value = 2.94
$ws.Cells.Item(2, 17).Value = 1.66
$ws.Cells.Item(2, 18).Value = 1.53
$ws.Cells.Item(2, 19).Value = 2.64

# Row 4: odds updates
$ws.Cells.Item(4, 31).Value = 1000

# Row 5: odds updates
$ws.Cells.Item(5, 20).Value = 2.1
$ws.Cells.Item(5, 21).Value = 1.75

# Row 6: odds updates
$ws.Cells.Item(6, 20).Value = 1.9

# Row 7: odds updates
$ws.Cells.Item(7, 6).Value = 5.8
$ws.Cells.Item(7, 7).Value = 8
$ws.Cells.Item(7, 9).Value = 1.68
$ws.Cells.Item(7, 11).Value = 5.4
$ws.Cells.Item(7, 12).Value = 1.01
$ws.Cells.Item(7, 13).Value = 1.04
$ws.Cells.Item(7, 14).Value = 4.2
$ws.Cells.Item(7, 15).Value = 1.25
$ws.Cells.Item(7, 17).Value = 1.64
$ws.Cells.Item(7, 18).Value = 1.43
$ws.Cells.Item(7, 19).Value = 2.58
$ws.Cells.Item(7, 20).Value = 1.8
$ws.Cells.Item(7, 21).Value = 2
$ws.Cells.Item(7, 22).Value = 2.46
$ws.Cells.Item(7, 23).Value = 1.16
$ws.Cells.Item(7, 24).Value = 19.5
$ws.Cells.Item(7, 25).Value = 9.6
$ws.Cells.Item(7, 26).Value = 10.5
$ws.Cells.Item(7, 27).Value = 16
$ws.Cells.Item(7, 28).Value = 24
$ws.Cells.Item(7, 29).Value = 10.5
$ws.Cells.Item(7, 30).Value = 10.5
$ws.Cells.Item(7, 31).Value = 17
$ws.Cells.Item(7, 32).Value = 65
$ws.Cells.Item(7, 33).Value = 25
$ws.Cells.Item(7, 34).Value = 22
$ws.Cells.Item(7, 35).Value = 48
$ws.Cells.Item(7, 36).Value = 1000
$ws.Cells.Item(7, 37).Value = 110
$ws.Cells.Item(7, 38).Value = 95
$ws.Cells.Item(7, 39).Value = 130
$ws.Cells.Item(7, 40).Value = 120
$ws.Cells.Item(7, 41).Value = 8.4

# Row 8: new match data
$ws.Cells.Item(8, 1).Value = "Serbian Super League"
$ws.Cells.Item(8, 2).Value = "2025-12-04"
$ws.Cells.Item(8, 3).Value = "15:00:00"
$ws.Cells.Item(8, 4).Value = "Crvena Zvezda"
$ws.Cells.Item(8, 5).Value = "Cukaricki"
$ws.Cells.Item(8, 6).Value = 1.14
$ws.Cells.Item(8, 7).Value = 1.36
$ws.Cells.Item(8, 8).Value = 9.199999999999999
$ws.Cells.Item(8, 9).Value = 1000
$ws.Cells.Item(8, 10).Value = 6.4
$ws.Cells.Item(8, 11).Value = 950
$ws.Cells.Item(8, 12).Value = 1.01
$ws.Cells.Item(8, 13).Value = 1.01
$ws.Cells.Item(8, 14).Value = 3
$ws.Cells.Item(8, 15).Value = 1.07
$ws.Cells.Item(8, 16).Value = 3
$ws.Cells.Item(8, 17).Value = 1.07
$ws.Cells.Item(8, 18).Value = 1.86
$ws.Cells.Item(8, 19).Value = 1.64
$ws.Cells.Item(8, 20).Value = 1.04
$ws.Cells.Item(8, 21).Value = 1.04
$ws.Cells.Item(8, 22).Value = 1.01
$ws.Cells.Item(8, 23).Value = 3.85
$ws.Cells.Item(8, 24).Value = 1000
$ws.Cells.Item(8, 25).Value = 1000
$ws.Cells.Item(8, 26).Value = 1000
$ws.Cells.Item(8, 27).Value = 1000
$ws.Cells.Item(8, 28).Value = 1000
$ws.Cells.Item(8, 29).Value = 1000
$ws.Cells.Item(8, 30).Value = 1000
$ws.Cells.Item(8, 31).Value = 1000
$ws.Cells.Item(8, 32).Value = 1000
$ws.Cells.Item(8, 33).Value = 1000
$ws.Cells.Item(8, 34).Value = 1000
$ws.Cells.Item(8, 35).Value = 1000
$ws.Cells.Item(8, 36).Value = 1000
$ws.Cells.Item(8, 37).Value = 1000
$ws.Cells.Item(8, 38).Value = 1000
$ws.Cells.Item(8, 39).Value = 1000
$ws.Cells.Item(8, 40).Value = 1000
$ws.Cells.Item(8, 41).Value = 1000

# Row 9: new match data
$ws.Cells.Item(9, 1).Value = "Italian Serie B"
$ws.Cells.Item(9, 2).Value = "2025-12-04"
$ws.Cells.Item(9, 3).Value = "15:30:00"
$ws.Cells.Item(9, 4).Value = "Juve Stabia"
$ws.Cells.Item(9, 5).Value = "SSD Bari"
$ws.Cells.Item(9, 6).Value = 1.94
$ws.Cells.Item(9, 7).Value = 2
$ws.Cells.Item(9, 8).Value = 4.8
$ws.Cells.Item(9, 9).Value = 5.3
$ws.Cells.Item(9, 10).Value = 3.15
$ws.Cells.Item(9, 11).Value = 3.45
$ws.Cells.Item(9, 12).Value = 1.01
$ws.Cells.Item(9, 13).Value = 1.1
$ws.Cells.Item(9, 14).Value = 2.7
$ws.Cells.Item(9, 15).Value = 1.5
$ws.Cells.Item(9, 16).Value = 1.58
$ws.Cells.Item(9, 17).Value = 2.48
$ws.Cells.Item(9, 18).Value = 1.21
$ws.Cells.Item(9, 19).Value = 5
$ws.Cells.Item(9, 20).Value = 1.94
$ws.Cells.Item(9, 21).Value = 1.62
$ws.Cells.Item(9, 22).Value = 1.23
$ws.Cells.Item(9, 23).Value = 2
$ws.Cells.Item(9, 24).Value = 9.4
$ws.Cells.Item(9, 25).Value = 13.5
$ws.Cells.Item(9, 26).Value = 38
$ws.Cells.Item(9, 27).Value = 180
$ws.Cells.Item(9, 28).Value = 6.8
$ws.Cells.Item(9, 29).Value = 7.6
$ws.Cells.Item(9, 30).Value = 21
$ws.Cells.Item(9, 31).Value = 1000
$ws.Cells.Item(9, 32).Value = 11
$ws.Cells.Item(9, 33).Value = 11.5
$ws.Cells.Item(9, 34).Value = 25
$ws.Cells.Item(9, 35).Value = 130
$ws.Cells.Item(9, 36).Value = 25
$ws.Cells.Item(9, 37).Value = 27
$ws.Cells.Item(9, 38).Value = 60
$ws.Cells.Item(9, 39).Value = 240
$ws.Cells.Item(9, 40).Value = 23
$ws.Cells.Item(9, 41).Value = 180

# Row 11: odds updates
$ws.Cells.Item(11, 19).Value = 2.44
$ws.Cells.Item(11, 27).Value = 340
$ws.Cells.Item(11, 36).Value = 13.5
$ws.Cells.Item(11, 38).Value = 27
$ws.Cells.Item(11, 39).Value = 90
$ws.Cells.Item(11, 40).Value = 5.5

# Row 13: odds updates
$ws.Cells.Item(13, 6).Value = 2.08
$ws.Cells.Item(13, 10).Value = 3.4
$ws.Cells.Item(13, 11).Value = 3.65
$ws.Cells.Item(13, 12).Value = 1.37
$ws.Cells.Item(13, 16).Value = 1.77
$ws.Cells.Item(13, 17).Value = 2.1
